$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change direction from "ltr" to "rtl" - this triggers recalculation of the
# dependent D6:D22 shared formula cells.
$ws.Range("C3").Value = "rtl"

# Update the selection to match the new active cell / selection range.
$ws.Range("K12").Select()
